# Slide 12, shape "TextBox 10" contains the JS code sample. In the third
# paragraph ("  let total = principal*Math.pow(1+interest,years);") the
# leading "  let " run gets merged into the following "total = principal*"
# run, which keeps its own run properties (dirty="0").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(5)

$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)

$leadRun = $para.Runs(1, 1)
$totalRun = $para.Runs(2, 1)

# Prepend the leading run's text onto the second run (preserving the
# second run's formatting/rPr), then remove the now-duplicated leading run.
$totalRun.Text = $leadRun.Text + $totalRun.Text
$leadRun.Text = ""
